# CompStat weekly refresh: new week's crime data collected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header strings: volume/issue number and the reporting week's date range.
# (The original shared strings are split across several runs that all share
# identical formatting, so a plain text replacement renders identically.)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/23/2024  Through  12/29/2024"

# ---------------------------------------------------------------------------
# Helper: force a cell to hold literal TEXT (not an auto-coerced number),
# then restore its original number-format/style by copying formats only
# from a nearby cell that already carries the desired style.
# ---------------------------------------------------------------------------
function Set-TextValue($addr, $text, $styleSource) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($styleSource).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

function Set-NumberValue($addr, $num, $styleSource) {
    $ws.Range($addr).Value = $num
    $ws.Range($styleSource).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# Plain numeric updates (style/type unchanged from the original workbook).
# ---------------------------------------------------------------------------
$ws.Range("N14").Value = -86.363636363636
$ws.Range("L15").Value = 5.555555555555
$ws.Range("M15").Value = 111.111111111111
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 186
$ws.Range("J16").Value = 198
$ws.Range("K16").Value = -6.060606060606
$ws.Range("L16").Value = 5.681818181818
$ws.Range("M16").Value = -50.267379679144
$ws.Range("N16").Value = -85.525291828793
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = -18.518518518518
$ws.Range("I17").Value = 308
$ws.Range("J17").Value = 315
$ws.Range("K17").Value = -2.222222222222
$ws.Range("L17").Value = 15.355805243445
$ws.Range("M17").Value = 48.792270531401
$ws.Range("N17").Value = -46.804835924006
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 303
$ws.Range("J18").Value = 281
$ws.Range("K18").Value = 7.829181494661
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = -35.940803382663
$ws.Range("N18").Value = -77.062831188493
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 65
$ws.Range("H19").Value = -40
$ws.Range("I19").Value = 690
$ws.Range("J19").Value = 792
$ws.Range("K19").Value = -12.878787878787
$ws.Range("L19").Value = 4.072398190045
$ws.Range("M19").Value = 37.176938369781
$ws.Range("N19").Value = 24.548736462093
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -44.444444444444
$ws.Range("I20").Value = 118
$ws.Range("J20").Value = 187
$ws.Range("K20").Value = -36.898395721925
$ws.Range("L20").Value = -35.519125683060
$ws.Range("M20").Value = -27.160493827160
$ws.Range("N20").Value = -87.406616862326
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -35.714285714285
$ws.Range("F21").Value = 93
$ws.Range("G21").Value = 135
$ws.Range("H21").Value = -31.111111111111
$ws.Range("I21").Value = 1627
$ws.Range("J21").Value = 1790
$ws.Range("K21").Value = -9.106145251396
$ws.Range("L21").Value = 1.244555071561
$ws.Range("M21").Value = -5.953757225433
$ws.Range("N21").Value = -65.617075232459
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -60
$ws.Range("J22").Value = 29
$ws.Range("K22").Value = -20.689655172413
$ws.Range("M22").Value = -47.727272727272
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = -41.176470588235
$ws.Range("I23").Value = 162
$ws.Range("J23").Value = 207
$ws.Range("K23").Value = -21.739130434782
$ws.Range("L23").Value = 2.531645569620
$ws.Range("M23").Value = 21.804511278195
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -11.111111111111
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = 34.482758620689
$ws.Range("I24").Value = 1117
$ws.Range("J24").Value = 1038
$ws.Range("K24").Value = 7.610789980732
$ws.Range("L24").Value = -6.213266162888
$ws.Range("M24").Value = -10.209003215434
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 316
$ws.Range("J25").Value = 166
$ws.Range("K25").Value = 90.361445783132
$ws.Range("L25").Value = 21.538461538461
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = -75
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = -8.571428571428
$ws.Range("I26").Value = 508
$ws.Range("J26").Value = 487
$ws.Range("K26").Value = 4.312114989733
$ws.Range("L26").Value = 3.462321792260
$ws.Range("M26").Value = 4.958677685950
$ws.Range("L27").Value = -11.538461538461
$ws.Range("C28").Value = 1
$ws.Range("I28").Value = 74
$ws.Range("K28").Value = 42.307692307692
$ws.Range("L28").Value = 60.869565217391
$ws.Range("L33").Value = 0

# ---------------------------------------------------------------------------
# Cells that flip from a numeric value to the "N/A" text placeholders
# ("0" / "***.*", shared strings already used elsewhere e.g. C14/E14).
# ---------------------------------------------------------------------------
Set-TextValue "G15" "0"     "C14"
Set-TextValue "H15" "***.*" "E14"
Set-TextValue "G27" "0"     "C14"
Set-TextValue "H27" "***.*" "E14"
Set-TextValue "D28" "0"     "C14"
Set-TextValue "E28" "***.*" "E14"

# ---------------------------------------------------------------------------
# Cells that flip from the "N/A" text placeholders back to real numbers.
# ---------------------------------------------------------------------------
Set-NumberValue "D22" 3    "C16"
Set-NumberValue "E22" -100 "E16"

Write-Host "CompStat weekly figures updated."
